$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: add a new "Michael is a device manager..." persona paragraph
# right after the "Tabitha is a data entry clerk..." paragraph and
# before the "Use Cases" heading paragraph.
# ---------------------------------------------------------------------
$tabithaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Tabitha is a data entry clerk*") {
        $tabithaIndex = $i
        break
    }
}

if ($tabithaIndex -ge 1) {
    $tabithaPara = $d.Paragraphs.Item($tabithaIndex)
    $tabithaPara.Range.InsertParagraphAfter()
    $michaelPara = $d.Paragraphs.Item($tabithaIndex + 1)
    $apos = [char]0x2019
    $michaelText = "Michael is a device manager, who is in charge of tracking device check-in and check-outs. He is also color blind and can" + $apos + "t see the color green."
    $michaelPara.Range.Text = $michaelText
}

# ---------------------------------------------------------------------
# Edit 2: add two new use-case entries at the end of the Use Cases
# section, right after the existing "System allows for cut/copy/paste
# functions in the memory editor view." use case and its following
# blank paragraph, before the document's final blank paragraph.
# ---------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*System allows for cut/copy/paste functions in the memory editor view*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ge 1) {
    # The blank paragraph immediately following the anchor use case.
    $blankIndex = $anchorIndex + 1
    $insertAfterPara = $d.Paragraphs.Item($blankIndex)

    $newLines = @(
        "",
        "Actor: User",
        "System Goal: Convert Four-Letter Words to Six-Letter Words",
        "Preconditions: User wants to be able change all lines of memory to six letters without doing it manually line by line in UVSim VM.",
        "System allows for automatic conversion with the click of a button.",
        "",
        "Actor: User",
        "System Goal: Save file as new file",
        "Preconditions: User wants to be able to save the current working file as a new file in a certain location with a certain name.",
        "System allows for SaveAs button that lets the user select the new file location and name."
    )

    $currentPara = $insertAfterPara
    $currentIndex = $blankIndex
    foreach ($line in $newLines) {
        $currentPara.Range.InsertParagraphAfter()
        $currentIndex = $currentIndex + 1
        $currentPara = $d.Paragraphs.Item($currentIndex)
        if ($line -ne "") {
            $currentPara.Range.Text = $line
        }
    }
}
